# Apply targeted updates to column F (dSF) on Sheet1, matching the diff:
# repull data, push all data, mean calculation

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value  = 0
$ws.Range("F11").Value = 3
$ws.Range("F20").Value = 2
$ws.Range("F21").Value = -3
$ws.Range("F22").Value = -4
$ws.Range("F24").Value = 4
$ws.Range("F31").Value = 3
$ws.Range("F34").Value = 0
$ws.Range("F38").Value = -4

$wb.Save()
